$p = $ppt.ActivePresentation

# Slide 7 - "Escena "B"" title placeholder was left empty; fill it in.
$titleRange7 = $p.Slides.Item(7).Shapes.Item(1).TextFrame.TextRange
$titleRange7.Text = "Escena “B”"
$titleRange7.LanguageID = "es-ES"

# Slide 8 - same title text as slide 7, also left empty.
$titleRange8 = $p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange
$titleRange8.Text = "Escena “B”"
$titleRange8.LanguageID = "es-ES"

# Slide 10 - "Fuentes" (sources) title placeholder was left empty.
$titleRange10 = $p.Slides.Item(10).Shapes.Item(1).TextFrame.TextRange
$titleRange10.Text = "Fuentes"
$titleRange10.LanguageID = "es-ES"
